$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Additional 11 test-result rows (396-406) appended to the "Scenario Results" sheet.
$data = @(
  @("Login with valid username and password", "PASSED", "chrome"),
  @("Login with valid username and password", "PASSED", "chrome"),
  @("Create a country", "PASSED", "chrome"),
  @("Create and Delete Cities", "PASSED", "chrome"),
  @("Create a country", "FAILED", "chrome"),
  @("Create and Delete Cities", "FAILED", "chrome"),
  @("Login with valid username and password", "FAILED", "chrome"),
  @("Login with valid username and password", "FAILED", "chrome"),
  @("Create a country", "FAILED", "chrome"),
  @("Create and Delete Cities", "FAILED", "chrome"),
  @("Login with valid username and password", "PASSED", "chrome")
)

$startRow = 396
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
